$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 25008.5
$ws.Range("I21").Value = 25008.5
$ws.Range("K21").Value = 25008.5
$ws.Range("M21").Value = -24540.5
$ws.Range("H23").Value = 25008.5
$ws.Range("I23").Value = 25008.5
$ws.Range("K23").Value = 25008.5
$ws.Range("M23").Value = -24774.5
$ws.Range("H32").Value = 5774.231
$ws.Range("J32").Value = 4809.3
$ws.Range("L32").Value = 4809.3
$ws.Range("N32").Value = -5461.3
$ws.Range("H62").Value = 5000.524
$ws.Range("I62").Value = 2769.1538
$ws.Range("J62").Value = 8626.5
$ws.Range("K62").Value = 2769.1538
$ws.Range("L62").Value = 8626.5
$ws.Range("M62").Value = -2145.1538
$ws.Range("N62").Value = -9874.5
$ws.Range("H65").Value = 5000.524
$ws.Range("I65").Value = 2769.1538
$ws.Range("J65").Value = 8626.5
$ws.Range("K65").Value = 13845.769
$ws.Range("L65").Value = 43132.5
$ws.Range("M65").Value = -10725.769
$ws.Range("N65").Value = -49372.5
$ws.Range("H113").Value = 4628.3335
$ws.Range("I113").Value = 1499
$ws.Range("J113").Value = 4912.8184
$ws.Range("K113").Value = 1499
$ws.Range("L113").Value = 4912.8184
$ws.Range("M113").Value = 1755
$ws.Range("N113").Value = -11420.8184
$ws.Range("H116").Value = 4553.0625
$ws.Range("I116").Value = 3486.9167
$ws.Range("J116").Value = 7751.5
$ws.Range("K116").Value = 3486.9167
$ws.Range("L116").Value = 7751.5
$ws.Range("M116").Value = -44.91670000000022
$ws.Range("N116").Value = -14635.5
$ws.Range("H125").Value = 2818.6667
$ws.Range("I125").Value = 2444
$ws.Range("K125").Value = 21996
$ws.Range("M125").Value = -19536
$ws.Range("H132").Value = 1286.8077
$ws.Range("I132").Value = 980.26086
$ws.Range("K132").Value = 2940.78258
$ws.Range("M132").Value = -410.7825800000001
$ws.Range("H138").Value = 2784.1123
$ws.Range("J138").Value = 2965.6165
$ws.Range("L138").Value = 8896.8495
$ws.Range("N138").Value = -19176.8495

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").Value = $null
$ws.Range("H32").Value = 5026.7837
$ws.Range("I32").Value = 3878.5151
$ws.Range("K32").Value = 3878.5151
$ws.Range("M32").Value = -3591.5151
$ws.Range("H34").Value = 259250
$ws.Range("J34").Value = 334000
$ws.Range("L34").Value = 334000
$ws.Range("N34").Value = -334542
$ws.Range("H74").Value = 19612734
$ws.Range("I74").Value = 33336462
$ws.Range("J74").Value = 7405.7144
$ws.Range("K74").Value = 33336462
$ws.Range("L74").Value = 7405.7144
$ws.Range("M74").Value = -33335588
$ws.Range("N74").Value = -9153.714400000001
$ws.Range("H77").Value = 19612734
$ws.Range("I77").Value = 33336462
$ws.Range("J77").Value = 7405.7144
$ws.Range("K77").Value = 166682310
$ws.Range("L77").Value = 37028.572
$ws.Range("M77").Value = -166677942
$ws.Range("N77").Value = -45764.572
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").Value = $null
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").Value = $null
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").Value = $null
$ws.Range("H122").Value = 3311.48
$ws.Range("I122").Value = 2787.4614
$ws.Range("K122").Value = 8362.3842
$ws.Range("M122").Value = -5912.3842
$ws.Range("H132").Value = 2363.6775
$ws.Range("I132").Value = 1698.6207
$ws.Range("K132").Value = 5095.8621
$ws.Range("M132").Value = -2565.8621

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3488.5
$ws.Range("I20").Value = 3350.2856
$ws.Range("J20").Value = 3972.25
$ws.Range("K20").Value = 3350.2856
$ws.Range("L20").Value = 3972.25
$ws.Range("M20").Value = -3103.2856
$ws.Range("N20").Value = -4466.25
$ws.Range("H86").Value = 5473.4443
$ws.Range("I86").Value = 3250.5715
$ws.Range("K86").Value = 3250.5715
$ws.Range("M86").Value = -2127.5715
$ws.Range("H89").Value = 5473.4443
$ws.Range("I89").Value = 3250.5715
$ws.Range("K89").Value = 16252.8575
$ws.Range("M89").Value = -10636.8575
$ws.Range("H105").Value = 22356.572
$ws.Range("I105").Value = 42498.6
$ws.Range("J105").Value = 11166.556
$ws.Range("K105").Value = 42498.6
$ws.Range("L105").Value = 11166.556
$ws.Range("M105").Value = -40751.6
$ws.Range("N105").Value = -14660.556
$ws.Range("H107").Value = 1794.8235
$ws.Range("I107").Value = 1536.5714
$ws.Range("K107").Value = 1536.5714
$ws.Range("M107").Value = 383.4286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 32583.105
$ws.Range("I31").Value = 4240.037
$ws.Range("J31").Value = 102152.45
$ws.Range("K31").Value = 4240.037
$ws.Range("L31").Value = 102152.45
$ws.Range("M31").Value = -3945.037
$ws.Range("N31").Value = -102742.45
$ws.Range("H34").Value = 32583.105
$ws.Range("I34").Value = 4240.037
$ws.Range("J34").Value = 102152.45
$ws.Range("K34").Value = 4240.037
$ws.Range("L34").Value = 102152.45
$ws.Range("M34").Value = -4038.037
$ws.Range("N34").Value = -102556.45
$ws.Range("H53").Value = 58467
$ws.Range("J53").Value = 58467
$ws.Range("L53").Value = 58467
$ws.Range("N53").Value = -59681
$ws.Range("H99").Value = 3841.2856
$ws.Range("I99").Value = 3600
$ws.Range("J99").Value = 3881.5
$ws.Range("K99").Value = 3600
$ws.Range("L99").Value = 3881.5
$ws.Range("M99").Value = -2102
$ws.Range("N99").Value = -6877.5
$ws.Range("H105").Value = 4256.75
$ws.Range("I105").Value = 1517.25
$ws.Range("K105").Value = 1517.25
$ws.Range("M105").Value = 229.75
$ws.Range("H118").Value = 49742
$ws.Range("J118").Value = 49742
$ws.Range("L118").Value = 49742
$ws.Range("N118").Value = -53056
$ws.Range("H126").Value = 3841.2856
$ws.Range("I126").Value = 3600
$ws.Range("J126").Value = 3881.5
$ws.Range("K126").Value = 10800
$ws.Range("L126").Value = 11644.5
$ws.Range("M126").Value = -8330
$ws.Range("N126").Value = -16584.5
$ws.Range("H132").Value = 3480.5938
$ws.Range("I132").Value = 2474.1365
$ws.Range("J132").Value = 5694.8
$ws.Range("K132").Value = 7422.4095
$ws.Range("L132").Value = 17084.4
$ws.Range("M132").Value = -4892.4095
$ws.Range("N132").Value = -22144.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 192.9
$ws.Range("J7").Value = 204.75
$ws.Range("L7").Value = 614.25
$ws.Range("N7").Value = -838.25
$ws.Range("H114").Value = 4949
$ws.Range("I114").Value = 950
$ws.Range("J114").Value = 6282
$ws.Range("K114").Value = 2850
$ws.Range("L114").Value = 18846
$ws.Range("M114").Value = 404
$ws.Range("N114").Value = -25354
$ws.Range("I122").Value = 1709911.4
$ws.Range("K122").Value = 15389202.6
$ws.Range("M122").Value = -15386752.6
$ws.Range("H129").Value = 5955825.5
$ws.Range("J129").Value = 10421875
$ws.Range("L129").Value = 31265625
$ws.Range("N129").Value = -31275625
$ws.Range("H131").Value = 8192731
$ws.Range("I131").Value = 31251118
$ws.Range("J131").Value = 5557486.5
$ws.Range("K131").Value = 93753354
$ws.Range("L131").Value = 16672459.5
$ws.Range("M131").Value = -93748314
$ws.Range("N131").Value = -16682539.5
$ws.Range("H132").Value = 4657.143
$ws.Range("I132").Value = 4183.1665
$ws.Range("K132").Value = 37648.4985
$ws.Range("M132").Value = -35118.4985
$ws.Range("H137").Value = 57247.555
$ws.Range("J137").Value = 102296.1
$ws.Range("L137").Value = 306888.3
$ws.Range("N137").Value = -317088.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 45829.832
$ws.Range("J15").Value = 45829.832
$ws.Range("L15").Value = 45829.832
$ws.Range("N15").Value = -46405.832
$ws.Range("H32").Value = 75000
$ws.Range("J32").Value = 75000
$ws.Range("L32").Value = 75000
$ws.Range("N32").Value = -75592
$ws.Range("H81").Value = 45829.832
$ws.Range("J81").Value = 45829.832
$ws.Range("L81").Value = 45829.832
$ws.Range("N81").Value = -47825.832
$ws.Range("H84").Value = 45829.832
$ws.Range("J84").Value = 45829.832
$ws.Range("L84").Value = 137489.496
$ws.Range("N84").Value = -147473.496
$ws.Range("H97").Value = 2142.3845
$ws.Range("I97").Value = 1784
$ws.Range("K97").Value = 1784
$ws.Range("M97").Value = -1288
$ws.Range("H99").Value = 35539.3
$ws.Range("I99").Value = 24476.5
$ws.Range("K99").Value = 24476.5
$ws.Range("M99").Value = -22230.5
$ws.Range("H102").Value = 2306.6296
$ws.Range("I102").Value = 1599.2821
$ws.Range("K102").Value = 1599.2821
$ws.Range("M102").Value = 22.7179000000001
$ws.Range("H107").Value = 401.72
$ws.Range("J107").Value = 395.22223
$ws.Range("L107").Value = 395.22223
$ws.Range("N107").Value = -4235.22223
$ws.Range("H113").Value = 2645.2083
$ws.Range("I113").Value = 2174
$ws.Range("J113").Value = 3304.9
$ws.Range("K113").Value = 2174
$ws.Range("L113").Value = 3304.9
$ws.Range("M113").Value = -4
$ws.Range("N113").Value = -7644.9
$ws.Range("H122").Value = 10795.571
$ws.Range("I122").Value = 14145
$ws.Range("K122").Value = 42435
$ws.Range("M122").Value = -39985

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 15872.066
$ws.Range("I40").Value = 20609.223
$ws.Range("J40").Value = 8766.333000000001
$ws.Range("K40").Value = 20609.223
$ws.Range("L40").Value = 8766.333000000001
$ws.Range("M40").Value = -20473.223
$ws.Range("N40").Value = -9038.333000000001
$ws.Range("H61").Value = 15925.637
$ws.Range("I61").Value = 19741.459
$ws.Range("K61").Value = 19741.459
$ws.Range("M61").Value = -19539.459
$ws.Range("H104").Value = 65177.5
$ws.Range("J104").Value = 65177.5
$ws.Range("L104").Value = 65177.5
$ws.Range("N104").Value = -72165.5
$ws.Range("H113").Value = 15925.637
$ws.Range("I113").Value = 19741.459
$ws.Range("K113").Value = 19741.459
$ws.Range("M113").Value = -17571.459
$ws.Range("H118").Value = 104409
$ws.Range("J118").Value = 104409
$ws.Range("L118").Value = 104409
$ws.Range("N118").Value = -107723
$ws.Range("H122").Value = 293061.94
$ws.Range("I122").Value = 1002237.5
$ws.Range("J122").Value = 9391.700000000001
$ws.Range("K122").Value = 3006712.5
$ws.Range("L122").Value = 28175.1
$ws.Range("M122").Value = -3004262.5
$ws.Range("N122").Value = -33075.10000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 10000
$ws.Range("I9").Value = 10000
$ws.Range("K9").Value = 10000
$ws.Range("M9").Value = -9860
$ws.Range("H107").Value = 538.3333
$ws.Range("I107").Value = 555.2941
$ws.Range("K107").Value = 1665.8823
$ws.Range("M107").Value = 254.1177000000002
$ws.Range("H126").Value = 1632.125
$ws.Range("I126").Value = 1420.9333
$ws.Range("J126").Value = 4800
$ws.Range("K126").Value = 4262.7999
$ws.Range("L126").Value = 14400
$ws.Range("M126").Value = -1792.7999
$ws.Range("N126").Value = -19340
$ws.Range("H132").Value = 9240.666999999999
$ws.Range("I132").Value = 4608.5
$ws.Range("K132").Value = 13825.5
$ws.Range("M132").Value = -11295.5
$ws.Range("H136").Value = 5599.8335
$ws.Range("I136").Value = 3399.25
$ws.Range("K136").Value = 10197.75
$ws.Range("M136").Value = -7647.75
